# Update F-column numeric values (counts/stats) across all four sheets
# per the commit "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 451
$ws.Range("F5").Value = 368
$ws.Range("F6").Value = 588
$ws.Range("F7").Value = 60
$ws.Range("F10").Value = 401
$ws.Range("F11").Value = 359
$ws.Range("F12").Value = 775
$ws.Range("F13").Value = 783
$ws.Range("F14").Value = 9
$ws.Range("F16").Value = 1544
$ws.Range("F17").Value = 1544
$ws.Range("F18").Value = 980
$ws.Range("F20").Value = 1366
$ws.Range("F22").Value = 363
$ws.Range("F25").Value = 114
$ws.Range("F26").Value = 6754
$ws.Range("F27").Value = 5161
$ws.Range("F28").Value = 8
$ws.Range("F29").Value = 151
$ws.Range("F32").Value = 217
$ws.Range("F37").Value = 1316
$ws.Range("F39").Value = 260
$ws.Range("F40").Value = 632
$ws.Range("F43").Value = 270
$ws.Range("F45").Value = 158
$ws.Range("F48").Value = 103
$ws.Range("F49").Value = 9

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 14
$ws.Range("F15").Value = 53

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 217
$ws.Range("F5").Value = 84

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 451
$ws.Range("F4").Value = 14
$ws.Range("F6").Value = 217
$ws.Range("F7").Value = 84
$ws.Range("F8").Value = 368
$ws.Range("F9").Value = 588
$ws.Range("F10").Value = 60
$ws.Range("F14").Value = 401
$ws.Range("F15").Value = 359
$ws.Range("F16").Value = 775
$ws.Range("F17").Value = 783
$ws.Range("F18").Value = 9
$ws.Range("F20").Value = 1544
$ws.Range("F21").Value = 1544
$ws.Range("F22").Value = 981
$ws.Range("F24").Value = 363
$ws.Range("F26").Value = 114
$ws.Range("F29").Value = 6754
$ws.Range("F30").Value = 5161
$ws.Range("F31").Value = 217
$ws.Range("F33").Value = 1316
$ws.Range("F36").Value = 260
$ws.Range("F38").Value = 632
$ws.Range("F41").Value = 53
$ws.Range("F43").Value = 270
$ws.Range("F44").Value = 158
$ws.Range("F47").Value = 103
$ws.Range("F50").Value = 9
